$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from an existing header cell
# onto the three new header cells before setting their text.
$headerStyleSource = $ws.Range("AC1")
$newHeaders = $ws.Range("AD1:AF1")
$headerStyleSource.Copy()
$newHeaders.PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every player row.
for ($r = 2; $r -le 60; $r++) {
    $ws.Cells.Item($r, 30).Value = 58
    $ws.Cells.Item($r, 31).Value = 104
    $ws.Cells.Item($r, 32).Value = 0
}
